# Atualização de bases das ligas, do dia: 31-03-2024 às 20:29
#
# For each listed pair of match rows, the recorded fixture id (column B) and
# every odds/result column from F through AC were swapped between the two
# rows (columns A/C/D/E - the running index, the two "Iraq League" markers,
# and the match date/time - stay put). A brand-new, not-yet-played fixture
# is appended as row 151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param(
        [int]$RowA,
        [int]$RowB
    )

    $idRangeA = $ws.Range("B$RowA")
    $idRangeB = $ws.Range("B$RowB")
    $dataRangeA = $ws.Range("F$RowA" + ":AC$RowA")
    $dataRangeB = $ws.Range("F$RowB" + ":AC$RowB")

    $idValueA = $idRangeA.Value2
    $idValueB = $idRangeB.Value2
    $dataValueA = $dataRangeA.Value2
    $dataValueB = $dataRangeB.Value2

    $idRangeA.Value2 = $idValueB
    $idRangeB.Value2 = $idValueA
    $dataRangeA.Value2 = $dataValueB
    $dataRangeB.Value2 = $dataValueA
}

Swap-MatchRows 17 18
Swap-MatchRows 22 23
Swap-MatchRows 41 42
Swap-MatchRows 54 55
Swap-MatchRows 73 74
Swap-MatchRows 78 79
Swap-MatchRows 103 104

# New fixture appended at the bottom of the table - the match has not been
# played yet, so there is no FTHG/FTAG/FTR (H/I/J) nor closing PL_Ahh/PL_Aha
# (AB/AC) value.
$newRow = 151

$ws.Range("A$newRow").Value2 = 149
$ws.Range("A$newRow").Borders.LineStyle = 1
$ws.Range("A$newRow").Font.Bold = $true
$ws.Range("A$newRow").HorizontalAlignment = -4108
$ws.Range("A$newRow").VerticalAlignment = -4160

$ws.Range("B$newRow").Value2 = 8030838
$ws.Range("C$newRow").Value2 = "Iraq League"
$ws.Range("D$newRow").Value2 = "Iraq League"

$ws.Range("E$newRow").Value2 = 45383.64583333334
$ws.Range("E$newRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("F$newRow").Value2 = "Al Naft SC"
$ws.Range("G$newRow").Value2 = "Al Kahrabaa"

$ws.Range("K$newRow").Value2 = 2.75
$ws.Range("L$newRow").Value2 = 2.875
$ws.Range("M$newRow").Value2 = 2.5
$ws.Range("N$newRow").Value2 = 2.75
$ws.Range("O$newRow").Value2 = 2.875
$ws.Range("P$newRow").Value2 = 2.5
$ws.Range("Q$newRow").Value2 = 0
$ws.Range("R$newRow").Value2 = 2
$ws.Range("S$newRow").Value2 = 1.8
$ws.Range("T$newRow").Value2 = 2
$ws.Range("U$newRow").Value2 = 1.85
$ws.Range("V$newRow").Value2 = 1.95
$ws.Range("W$newRow").Value2 = 0
$ws.Range("X$newRow").Value2 = 0
$ws.Range("Y$newRow").Value2 = 0
$ws.Range("Z$newRow").Value2 = 0
$ws.Range("AA$newRow").Value2 = 0
